$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.602.30"
$ws.Range("E2").Value = "  +0.41%  "

$ws.Range("D3").Value = "1.811.43"
$ws.Range("E3").Value = "  +0.39%  "

$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.61"
$ws.Range("E5").Value = "  -1.21%  "

$ws.Range("E6").Value = "  +3.15%  "

$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "37.32"
$ws.Range("E8").Value = "  +7.13%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.293"
$ws.Range("E9").Value = "  -2.31%  "

$ws.Range("E10").Value = "  -1.58%  "

$ws.Range("E11").Value = "  +1.45%  "

$ws.Range("D12").Value = "2.074.63"
$ws.Range("E12").Value = "  +0.45%  "

$ws.Range("E13").Value = "  +1.79%  "

$ws.Range("D14").Value = "1.806.73"
$ws.Range("E14").Value = "  +0.16%  "

$ws.Range("E15").Value = "  -1.30%  "

$ws.Range("D16").Value = "34.550.16"
$ws.Range("E16").Value = "  +0.29%  "

$ws.Range("E17").Value = "  +1.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.70"
$ws.Range("E18").Value = "  -0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.26"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("E20").Value = "  -2.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.25"
$ws.Range("E21").Value = "  -1.81%  "

$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.13"
$ws.Range("E23").Value = "  -0.90%  "

$ws.Range("E24").Value = "  +4.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.09"
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.29"
$ws.Range("E27").Value = "  +2.86%  "

$ws.Range("E28").Value = "  +2.18%  "

$ws.Range("E29").Value = "  -0.24%  "

$ws.Range("E30").Value = "  -0.22%  "

$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0519"
$ws.Range("E33").Value = "  -2.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -0.40%  "

$ws.Range("D35").Value = "1.367.10"
$ws.Range("E35").Value = "  -1.94%  "

$ws.Range("E36").Value = "  -4.30%  "

$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("E38").Value = "  -4.43%  "

$ws.Range("E39").Value = "  -1.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.42"
$ws.Range("E40").Value = "  +1.24%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "81.19"
$ws.Range("E41").Value = "  -2.74%  "

$ws.Range("E42").Value = "  -1.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.940"
$ws.Range("E43").Value = "  -1.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.18"
$ws.Range("E44").Value = "  +5.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.82"
$ws.Range("E45").Value = "  +1.88%  "

$ws.Range("E46").Value = "  -1.97%  "

$ws.Range("D47").Value = "1.974.02"
$ws.Range("E47").Value = "  +0.44%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.82"
$ws.Range("E48").Value = "  -3.41%  "

$ws.Range("E49").Value = "  -0.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "102.69"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("E51").Value = "  -7.79%  "
